# feat: add 2022-Q1 data
#
# - Inserts a "2022-Q1" sheet (fund holdings) between "2021-Q4" and "总计".
# - Updates "总计" with a new summary row for 2022-Q1 (now listed first,
#   2021-Q4 shifts down to row 3).
#
# Implementation note: worksheet handles returned by Worksheets.Item(...) in
# this host are positional (re-resolved by current index), not bound to a
# stable object identity. Any handle captured before a sheet-order-changing
# operation (Add/rename-that-reorders/etc.) becomes stale, so sheets are
# re-fetched **by name** after such operations instead of being cached.
#
# Implementation note 2: sheetId is assigned sequentially
# (max-existing-id + 1) by Worksheets.Add, regardless of where the sheet is
# inserted positionally. To reproduce the target sheetId ordering
# (2021-Q4=1, 2022-Q1=2, 总计=3) the existing "总计" sheet (sheetId=2) is
# renamed in place to become "2022-Q1" and repopulated, then a brand new
# sheet is appended at the end and named "总计" (picking up sheetId=3).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Turn the existing "总计" sheet into "2022-Q1" (keeps sheetId=2, plus all
#    of its existing sheetPr/pageMargins boilerplate) and populate it with
#    the fund holdings table.
# ---------------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item("总计")
$wsQ1.Name = "2022-Q1"

$wsQ4 = $wb.Worksheets.Item("2021-Q4")

# Header row
$wsQ1.Cells.Item(1,2).Value = "基金代码"
$wsQ1.Cells.Item(1,3).Value = "基金名称"
$wsQ1.Cells.Item(1,4).Value = "基金规模"
$wsQ1.Cells.Item(1,5).Value = "股票总仓位"
$wsQ1.Cells.Item(1,6).Value = "仓位占比"
$wsQ1.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ1.Cells.Item(1,8).Value = "仓位排名"

# Copy the header style (bold + border) from the "2021-Q4" sheet's header row
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns B (fund code) and D:G (numeric-looking ratios) must be stored as
# TEXT, matching how the sibling quarter sheet stores them (and preserving
# leading zeros in fund codes) - pre-format as Text before writing values.
$wsQ1.Range("B2:B4").NumberFormat = "@"
$wsQ1.Range("D2:G4").NumberFormat = "@"

# Copy column-A's numeric index style down onto rows 3:4 (row 2 already
# carries it from the original "总计" sheet)
$wsQ1.Range("A2").Copy()
$wsQ1.Range("A3:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$fundRows = @(
    @(0, "006233", "万家汽车新趋势混合A", "4.33", "92.45", "2.46", "0.1065", 8),
    @(1, "007251", "广发睿享稳健增利混合", "3.69", "38.80", "2.26", "0.0834", 2),
    @(2, "006234", "万家汽车新趋势混合C", "2.52", "92.45", "2.46", "0.0620", 8)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $data = $fundRows[$i]
    $wsQ1.Cells.Item($r,1).Value = $data[0]
    $wsQ1.Cells.Item($r,2).Value = $data[1]
    $wsQ1.Cells.Item($r,3).Value = $data[2]
    $wsQ1.Cells.Item($r,4).Value = $data[3]
    $wsQ1.Cells.Item($r,5).Value = $data[4]
    $wsQ1.Cells.Item($r,6).Value = $data[5]
    $wsQ1.Cells.Item($r,7).Value = $data[6]
    $wsQ1.Cells.Item($r,8).Value = $data[7]
}

# ---------------------------------------------------------------------------
# 2) Append a brand-new "总计" sheet (picks up the next sheetId) summarizing
#    both quarters, with 2022-Q1 listed first and 2021-Q4 shifted to row 3.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsTotal = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsTotal.Name = "总计"

# Match page margins used by the rest of the workbook (0.75/0.75/1/1/.5/.5 in)
$wsTotal.PageSetup.LeftMargin   = 54
$wsTotal.PageSetup.RightMargin  = 54
$wsTotal.PageSetup.TopMargin    = 72
$wsTotal.PageSetup.BottomMargin = 72
$wsTotal.PageSetup.HeaderMargin = 36
$wsTotal.PageSetup.FooterMargin = 36

$wsQ4 = $wb.Worksheets.Item("2021-Q4")

# Header row (copy text + style from the old "总计" header, now "2021-Q4"
# style source)
$wsTotal.Cells.Item(1,2).Value = "日期"
$wsTotal.Cells.Item(1,3).Value = "持有数量(只)"
$wsTotal.Cells.Item(1,4).Value = "持有市值(亿元)"
$wsQ4.Range("B1:D1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column-A numeric index style
$wsQ4.Range("A2").Copy()
$wsTotal.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2022-Q1 summary (new, now first)
$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q1"
$wsTotal.Cells.Item(2,3).Value = 3
$wsTotal.Cells.Item(2,4).Value = 0.25

# 2021-Q4 summary (existing, shifted to row 3)
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(3,2).Value = "2021-Q4"
$wsTotal.Cells.Item(3,3).Value = 5
$wsTotal.Cells.Item(3,4).Value = 1.1
